# Auto-generated edit script applying cached market-price recalculation
# updates across the per-job Leve profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Each H:N cell holds a literal cached value (no formulas in this workbook), so every
# changed cell is written directly; cells that disappear/appear in the diff are cleared
# or newly populated to match.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value2 = 135.85715
$ws.Range("I4").Value2 = 135.85715
$ws.Range("K4").Value2 = 135.85715
$ws.Range("M4").Value2 = -21.85714999999999
$ws.Range("H52").Value2 = 2520.7144
$ws.Range("I52").Value2 = 3022.3333
$ws.Range("J52").Value2 = 2144.5
$ws.Range("K52").Value2 = 9066.999899999999
$ws.Range("L52").Value2 = 6433.5
$ws.Range("M52").Value2 = -8906.999899999999
$ws.Range("N52").Value2 = -6753.5
$ws.Range("H53").Value2 = 320.53845
$ws.Range("I53").Value2 = 226
$ws.Range("J53").Value2 = 430.83334
$ws.Range("K53").Value2 = 226
$ws.Range("L53").Value2 = 430.83334
$ws.Range("M53").Value2 = 411
$ws.Range("N53").Value2 = -1704.83334
$ws.Range("H70").Value2 = 102307.8
$ws.Range("J70").Value2 = 201948.2
$ws.Range("L70").Value2 = 605844.6000000001
$ws.Range("N70").Value2 = -606384.6000000001
$ws.Range("H73").Value2 = 102307.8
$ws.Range("J73").Value2 = 201948.2
$ws.Range("L73").Value2 = 605844.6000000001
$ws.Range("N73").Value2 = -607716.6000000001
$ws.Range("H88").Value2 = 2854.4666
$ws.Range("I88").Value2 = 1399.8334
$ws.Range("J88").Value2 = 3824.2222
$ws.Range("K88").Value2 = 1399.8334
$ws.Range("L88").Value2 = 3824.2222
$ws.Range("M88").Value2 = -993.8334
$ws.Range("N88").Value2 = -4636.2222
$ws.Range("H91").Value2 = 2854.4666
$ws.Range("I91").Value2 = 1399.8334
$ws.Range("J91").Value2 = 3824.2222
$ws.Range("K91").Value2 = 1399.8334
$ws.Range("L91").Value2 = 3824.2222
$ws.Range("M91").Value2 = 4.166600000000017
$ws.Range("N91").Value2 = -6632.2222
$ws.Range("H107").Value2 = 2294.5
$ws.Range("I107").Value2 = 2059.6667
$ws.Range("J107").Value2 = 2999
$ws.Range("K107").Value2 = 2059.6667
$ws.Range("L107").Value2 = 2999
$ws.Range("M107").Value2 = -139.6667000000002
$ws.Range("N107").Value2 = -6839
$ws.Range("H132").Value2 = 3174.9783
$ws.Range("I132").Value2 = 3266.0715
$ws.Range("J132").Value2 = 2218.5
$ws.Range("K132").Value2 = 9798.2145
$ws.Range("L132").Value2 = 6655.5
$ws.Range("M132").Value2 = -7268.2145
$ws.Range("N132").Value2 = -11715.5
$ws.Range("H137").Value2 = 1101.6666
$ws.Range("I137").Value2 = 1102.5714
$ws.Range("J137").Value2 = 1098.5
$ws.Range("K137").Value2 = 3307.7142
$ws.Range("L137").Value2 = 3295.5
$ws.Range("M137").Value2 = -757.7142000000003
$ws.Range("N137").Value2 = -8395.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value2 = 2470.65
$ws.Range("I32").Value2 = 784.4865
$ws.Range("J32").Value2 = 23266.666
$ws.Range("K32").Value2 = 784.4865
$ws.Range("L32").Value2 = 23266.666
$ws.Range("M32").Value2 = -497.4865
$ws.Range("N32").Value2 = -23840.666
$ws.Range("H61").Value2 = 8848.5
$ws.Range("I61").Value2 = 8848.5
$ws.Range("K61").Value2 = 8848.5
$ws.Range("M61").Value2 = -8636.5
$ws.Range("H63").Value2 = 6197.3335
$ws.Range("I63").Value2 = 5442.591
$ws.Range("K63").Value2 = 5442.591
$ws.Range("M63").Value2 = -4756.591
$ws.Range("H66").Value2 = 6197.3335
$ws.Range("I66").Value2 = 5442.591
$ws.Range("K66").Value2 = 27212.955
$ws.Range("M66").Value2 = -23780.955
$ws.Range("H124").Value2 = 39249.25
$ws.Range("J124").Value2 = 39249.25
$ws.Range("L124").Value2 = 39249.25
$ws.Range("N124").Value2 = -49069.25
$ws.Range("H136").Value2 = 8848.5
$ws.Range("I136").Value2 = 8848.5
$ws.Range("K136").Value2 = 26545.5
$ws.Range("M136").Value2 = -23995.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value2 = 1761.8695
$ws.Range("I107").Value2 = 1773.7727
$ws.Range("K107").Value2 = 1773.7727
$ws.Range("M107").Value2 = 146.2273

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value2 = 9569.235000000001
$ws.Range("I31").Value2 = 10491.8
$ws.Range("K31").Value2 = 10491.8
$ws.Range("M31").Value2 = -10196.8
$ws.Range("H34").Value2 = 9569.235000000001
$ws.Range("I34").Value2 = 10491.8
$ws.Range("K34").Value2 = 10491.8
$ws.Range("M34").Value2 = -10289.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value2 = 902.6316
$ws.Range("I2").Value2 = 1055.3334
$ws.Range("J2").Value2 = 330
$ws.Range("K2").Value2 = 6332.0004
$ws.Range("L2").Value2 = 1980
$ws.Range("M2").Value2 = -6219.0004
$ws.Range("N2").Value2 = -2206
$ws.Range("H4").Value2 = 58376550
$ws.Range("I4").Value2 = 125086500
$ws.Range("J4").Value2 = 5008594.5
$ws.Range("K4").Value2 = 375259500
$ws.Range("L4").Value2 = 15025783.5
$ws.Range("M4").Value2 = -375259388
$ws.Range("N4").Value2 = -15026007.5
$ws.Range("H9").Value2 = 15499.333
$ws.Range("J9").Value2 = 10249.5
$ws.Range("L9").Value2 = 30748.5
$ws.Range("N9").Value2 = -31196.5
$ws.Range("H10").Value2 = 419.6
$ws.Range("I10").Value2 = 149.5
$ws.Range("K10").Value2 = 448.5
$ws.Range("M10").Value2 = -309.5
$ws.Range("H15").Value2 = 3693.6667
$ws.Range("I15").Value2 = 4999
$ws.Range("J15").Value2 = 3636.913
$ws.Range("K15").Value2 = 14997
$ws.Range("L15").Value2 = 10910.739
$ws.Range("M15").Value2 = -14857
$ws.Range("N15").Value2 = -11190.739
$ws.Range("H16").Value2 = 7849.75
$ws.Range("I16").Value2 = 1350
$ws.Range("J16").Value2 = 14349.5
$ws.Range("K16").Value2 = 4050
$ws.Range("L16").Value2 = 43048.5
$ws.Range("M16").Value2 = -3877
$ws.Range("N16").Value2 = -43394.5
$ws.Range("H20").Value2 = 10000
$ws.Range("I20").Value2 = 10000
$ws.Range("K20").Value2 = 30000
$ws.Range("M20").Value2 = -29773
$ws.Range("H21").Value2 = 9592.333000000001
$ws.Range("I21").Value2 = 78
$ws.Range("J21").Value2 = 14349.5
$ws.Range("K21").Value2 = 234
$ws.Range("L21").Value2 = 43048.5
$ws.Range("M21").Value2 = -61
$ws.Range("N21").Value2 = -43394.5
$ws.Range("H22").Value2 = 5000
$ws.Range("J22").Value2 = 5000
$ws.Range("L22").Value2 = 15000
$ws.Range("N22").Value2 = -15338
$ws.Range("H26").Value2 = 2035.7778
$ws.Range("I26").Value2 = 42
$ws.Range("J26").Value2 = 2605.4285
$ws.Range("K26").Value2 = 126
$ws.Range("L26").Value2 = 7816.2855
$ws.Range("M26").Value2 = 162
$ws.Range("N26").Value2 = -8392.2855
$ws.Range("H27").Value2 = 5000
$ws.Range("J27").Value2 = 5000
$ws.Range("L27").Value2 = 15000
$ws.Range("N27").Value2 = -15204
$ws.Range("H33").Value2 = 102.25
$ws.Range("J33").Value2 = 0
$ws.Range("L33").Value2 = 0
$ws.Range("N33").ClearContents()
$ws.Range("H34").Value2 = 6493.3
$ws.Range("J34").Value2 = 7103.778
$ws.Range("L34").Value2 = 21311.334
$ws.Range("N34").Value2 = -21479.334
$ws.Range("H38").Value2 = 142.75
$ws.Range("I38").Value2 = 84.2
$ws.Range("J38").Value2 = 169.36363
$ws.Range("K38").Value2 = 252.6
$ws.Range("L38").Value2 = 508.09089
$ws.Range("M38").Value2 = 94.39999999999998
$ws.Range("N38").Value2 = -1202.09089
$ws.Range("H39").Value2 = 7036.375
$ws.Range("I39").Value2 = 2297.5
$ws.Range("J39").Value2 = 8616
$ws.Range("K39").Value2 = 6892.5
$ws.Range("L39").Value2 = 25848
$ws.Range("M39").Value2 = -6598.5
$ws.Range("N39").Value2 = -26436
$ws.Range("H40").Value2 = 179.28572
$ws.Range("I40").Value2 = 63.75
$ws.Range("K40").Value2 = 255
$ws.Range("M40").Value2 = -186
$ws.Range("H98").Value2 = 1413.9166
$ws.Range("I98").Value2 = 1188.1666
$ws.Range("K98").Value2 = 3564.4998
$ws.Range("M98").Value2 = -2066.4998
$ws.Range("H132").Value2 = 4999
$ws.Range("J132").Value2 = 0
$ws.Range("L132").Value2 = 0
$ws.Range("N132").ClearContents()
$ws.Range("H138").Value2 = 17299.111
$ws.Range("J138").Value2 = 19956
$ws.Range("L138").Value2 = 59868
$ws.Range("N138").Value2 = -70148

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value2 = 5598.4287
$ws.Range("I70").Value2 = 5448.1665
$ws.Range("K70").Value2 = 5448.1665
$ws.Range("M70").Value2 = -5178.1665
$ws.Range("H73").Value2 = 5598.4287
$ws.Range("I73").Value2 = 5448.1665
$ws.Range("K73").Value2 = 5448.1665
$ws.Range("M73").Value2 = -4512.1665
$ws.Range("H102").Value2 = 3271.4546
$ws.Range("I102").Value2 = 3271.4546
$ws.Range("K102").Value2 = 3271.4546
$ws.Range("M102").Value2 = -1649.4546
$ws.Range("H132").Value2 = 3346.9285
$ws.Range("I132").Value2 = 3419.2
$ws.Range("J132").Value2 = 2744.6667
$ws.Range("K132").Value2 = 10257.6
$ws.Range("L132").Value2 = 8234.000100000001
$ws.Range("M132").Value2 = -7727.599999999999
$ws.Range("N132").Value2 = -13294.0001
$ws.Range("H137").Value2 = 88952.39999999999
$ws.Range("J137").Value2 = 88952.39999999999
$ws.Range("L137").Value2 = 88952.39999999999
$ws.Range("N137").Value2 = -99152.39999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value2 = 5203.625
$ws.Range("I68").Value2 = 5203.625
$ws.Range("J68").Value2 = 0
$ws.Range("K68").Value2 = 5203.625
$ws.Range("L68").Value2 = 0
$ws.Range("M68").Value2 = -4454.625
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value2 = 5203.625
$ws.Range("I71").Value2 = 5203.625
$ws.Range("J71").Value2 = 0
$ws.Range("K71").Value2 = 26018.125
$ws.Range("L71").Value2 = 0
$ws.Range("M71").Value2 = -22274.125
$ws.Range("N71").ClearContents()
$ws.Range("H93").Value2 = 2712.6667
$ws.Range("J93").Value2 = 5056
$ws.Range("L93").Value2 = 5056
$ws.Range("N93").Value2 = -7552

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value2 = 4872.5
$ws.Range("I62").Value2 = 4872.5
$ws.Range("J62").Value2 = 0
$ws.Range("K62").Value2 = 4872.5
$ws.Range("L62").Value2 = 0
$ws.Range("M62").Value2 = -4248.5
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value2 = 4872.5
$ws.Range("I65").Value2 = 4872.5
$ws.Range("J65").Value2 = 0
$ws.Range("K65").Value2 = 24362.5
$ws.Range("L65").Value2 = 0
$ws.Range("M65").Value2 = -21242.5
$ws.Range("N65").ClearContents()
$ws.Range("H100").Value2 = 850
$ws.Range("I100").Value2 = 1400
$ws.Range("J100").Value2 = 300
$ws.Range("K100").Value2 = 2800
$ws.Range("L100").Value2 = 600
$ws.Range("M100").Value2 = -2259
$ws.Range("N100").Value2 = -1682
